# Updates license text from "CC BY" to "CC BY SA" on both footer shapes
# (slide 1 shape id=143, slide 2 shape id=261), matching the commit:
# "Updates license on all sheets from CC BY to CC BY SA"

$p = $ppt.ActivePresentation

function Find-ShapeById($slide, $targetId) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $targetId) {
            return $candidate
        }
    }
    return $null
}

# --- Slide 1 / shape id 143 : "CC BY " + "RStudio" -> "CC BY SA" + " RStudio" ---
$slide1 = $p.Slides.Item(1)
$shape143 = Find-ShapeById $slide1 143

$shape143.Name = "RStudio® is a trademark of RStudio, Inc.  •  CC BY SA RStudio •  info@rstudio.com  •  844-448-1212 • rstudio.com •  Learn more with browseVignettes(package = c(`"dplyr`", `"tibble`"))  •  dplyr  0.5.0 •  tibble  1.2.0  •  Updated: 2017-01"

# the textbox has spAutoFit; stash geometry so the re-layout triggered by
# the run-text edits below doesn't drift the shape's stored height/width
$origHeight143 = $shape143.Height
$origWidth143 = $shape143.Width

$tr143 = $shape143.TextFrame.TextRange
$tr143.Runs(2).Text = "CC BY SA"
$tr143.Runs(3).Text = " RStudio •  "

$shape143.Height = $origHeight143
$shape143.Width = $origWidth143

# --- Slide 2 / shape id 261 : "CC BY " + "RStudio" -> "CC BY SA" + "  RStudio" ---
$slide2 = $p.Slides.Item(2)
$shape261 = Find-ShapeById $slide2 261

$shape261.Name = "RStudio® is a trademark of RStudio, Inc.  •  CC BY SA  RStudio •  info@rstudio.com  •  844-448-1212 • rstudio.com •  Learn more with browseVignettes(package = c(`"dplyr`", `"tibble`"))  •  dplyr  0.5.0 •  tibble  1.2.0  •  Updated: 2017-01"

$origHeight261 = $shape261.Height
$origWidth261 = $shape261.Width

$tr261 = $shape261.TextFrame.TextRange
$tr261.Runs(2).Text = "CC BY SA"
$tr261.Runs(3).Text = "  RStudio •  "

$shape261.Height = $origHeight261
$shape261.Width = $origWidth261
